# Generate Report for Handoff
# Updates the localization status report: the item has moved from
# "In Translation" to "Ready for handoff", and the handoff timestamps
# are refreshed accordingly on every sheet that tracks them.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-26 00:58:12"

# --- zh-cn sheet ------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-26 00:58:07"

# --- de-de sheet ------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-26 00:58:12"

# --- Re-fit the "Status"/language columns now that the new text is
#     wider than before ("Ready for handoff" vs "In Translation") ----
$newColumnWidth = 16.333333333333332
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
